$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows before row 51 (shifts old 51,52 down to 53,54)
$ws.Rows.Item(51).Insert()
$ws.Rows.Item(51).Insert()

# Update D51/E51 (was 3.2E-2, now 0.02)
$ws.Range("D51").Value = 0.02
$ws.Range("E51").Value = 0.72213662044170501

# New row 52
$ws.Range("D52").Value = 0.025000000000000001
$ws.Range("E52").Value = 0.72213662044170501

# Update selection to match target diff
$ws.Range("E50:E52").Select()
